$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E10").Value = "ignore"
$ws.Range("E10").Select()
